$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.868.37'
$ws.Range("E2").Value = '  -1.15%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.428.51'
$ws.Range("E3").Value = '  -1.58%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '574.15'
$ws.Range("E5").Value = '  -0.90%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.06'
$ws.Range("E6").Value = '  -0.92%  '
$ws.Range("E7").Value = '  +3.54%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.429.57'
$ws.Range("E9").Value = '  -1.60%  '
$ws.Range("E10").Value = '  -1.34%  '
$ws.Range("E11").Value = '  -1.85%  '
$ws.Range("E12").Value = '  -0.60%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.017.70'
$ws.Range("E13").Value = '  -1.64%  '
$ws.Range("E14").Value = '  -0.05%  '
$ws.Range("E15").Value = '  -3.14%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.73'
$ws.Range("E16").Value = '  -3.60%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.836.57'
$ws.Range("E17").Value = '  -1.23%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.431.72'
$ws.Range("E18").Value = '  -0.99%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.37'
$ws.Range("E19").Value = '  -0.82%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.88'
$ws.Range("E20").Value = '  -2.92%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '381.39'
$ws.Range("E21").Value = '  -2.41%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.02'
$ws.Range("E22").Value = '  -2.69%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.549'
$ws.Range("E23").Value = '  -0.16%  '
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.39'
$ws.Range("E25").Value = '  -1.68%  '
$ws.Range("E26").Value = '  -3.79%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.13'
$ws.Range("E27").Value = '  +5.38%  '
$ws.Range("E28").Value = '  -0.45%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("E30").Value = '  +3.15%  '
$ws.Range("E31").Value = '  -3.36%  '
$ws.Range("E32").Value = '  -2.31%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.29'
$ws.Range("E33").Value = '  -1.74%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.09'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.61'
$ws.Range("E35").Value = '  +2.91%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '160.24'
$ws.Range("E36").Value = '  -1.65%  '
$ws.Range("E37").Value = '  -2.49%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.931.94'
$ws.Range("E38").Value = '  -4.61%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0759'
$ws.Range("E39").Value = '  -1.53%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.76'
$ws.Range("E40").Value = '  +4.24%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '26.44'
$ws.Range("E41").Value = '  -2.74%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.60'
$ws.Range("E42").Value = '  +1.64%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '43.00'
$ws.Range("E43").Value = '  +0.65%  '
$ws.Range("E44").Value = '  -1.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.772'
$ws.Range("E45").Value = '  -0.62%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '25.94'
$ws.Range("E46").Value = '  +0.24%  '
$ws.Range("E47").Value = '  +1.87%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '317.79'
$ws.Range("E48").Value = '  +1.77%  '
$ws.Range("E49").Value = '  -4.30%  '
$ws.Range("E50").Value = '  +0.77%  '
$ws.Range("E51").Value = '  -2.09%  '
